$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for "Terminal La Palmera de La Serena - Arveja Verde" need the
# values of row 2 and row 3 swapped for columns D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado) and
# P (Precio $/Kg), reflecting the weekly data being reordered.

$columns = @("D", "J", "K", "L", "M", "P")

foreach ($col in $columns) {
    $addr2 = "${col}2"
    $addr3 = "${col}3"
    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}
